# Automatic update of files.
# Rotate the A/Q/R/Y/AA values among rows 2-5: each row takes on the
# values that originally belonged to the next row (row 5 wraps to row 2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the original values for the cells that move, before any writes.
# Use .Value2 (not .Value) -- it returns plain scalars from this host.
$rows = 2, 3, 4, 5
$orig = @{}
foreach ($r in $rows) {
    $orig[$r] = @{
        A  = $ws.Range("A$r").Value2
        Q  = $ws.Range("Q$r").Value2
        R  = $ws.Range("R$r").Value2
        Y  = $ws.Range("Y$r").Value2
        AA = $ws.Range("AA$r").Value2
    }
}

# Map each destination row to the row whose original values it should receive.
$srcRow = @{ 2 = 3; 3 = 4; 4 = 5; 5 = 2 }

foreach ($r in $rows) {
    $src = $orig[$srcRow[$r]]

    $ws.Range("A$r").Value2 = $src.A
    $ws.Range("Q$r").Value2 = $src.Q
    $ws.Range("R$r").Value2 = $src.R

    # Force these as plain text so Excel doesn't reinterpret the
    # yyyy-mm-dd-looking string as a date serial number.
    $ws.Range("Y$r").NumberFormat = "@"
    $ws.Range("Y$r").Value2 = $src.Y
    $ws.Range("Y$r").Style = "Normal"

    $ws.Range("AA$r").NumberFormat = "@"
    $ws.Range("AA$r").Value2 = $src.AA
    $ws.Range("AA$r").Style = "Normal"
}
